$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value2 = 443.46155
$ws.Range("I19").Value2 = 459
$ws.Range("J19").Value2 = 433.75
$ws.Range("K19").Value2 = 459
$ws.Range("L19").Value2 = 433.75
$ws.Range("M19").Value2 = -284
$ws.Range("N19").Value2 = -783.75
$ws.Range("H28").Value2 = 1702.2667
$ws.Range("I28").Value2 = 428
$ws.Range("J28").Value2 = 9985
$ws.Range("K28").Value2 = 428
$ws.Range("L28").Value2 = 9985
$ws.Range("M28").Value2 = 57
$ws.Range("N28").Value2 = -10955
$ws.Range("H111").Value2 = 2000
$ws.Range("I111").Value2 = 2000
$ws.Range("J111").Value2 = 2000
$ws.Range("K111").Value2 = 6000
$ws.Range("L111").Value2 = 6000
$ws.Range("M111").Value2 = -2933
$ws.Range("N111").Value2 = -12134
$ws.Range("H135").Value2 = 1236.1
$ws.Range("I135").Value2 = 1370.5
$ws.Range("J135").Value2 = 1146.5
$ws.Range("K135").Value2 = 12334.5
$ws.Range("L135").Value2 = 10318.5
$ws.Range("M135").Value2 = -9799.5
$ws.Range("N135").Value2 = -15388.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value2 = 160.5
$ws.Range("I5").Value2 = 480
$ws.Range("J5").Value2 = 54
$ws.Range("K5").Value2 = 480
$ws.Range("L5").Value2 = 54
$ws.Range("M5").Value2 = -368
$ws.Range("N5").Value2 = -278
$ws.Range("H45").Value2 = 1114.125
$ws.Range("I45").Value2 = 967.6667
$ws.Range("J45").Value2 = 1553.5
$ws.Range("K45").Value2 = 967.6667
$ws.Range("L45").Value2 = 1553.5
$ws.Range("M45").Value2 = -590.6667
$ws.Range("N45").Value2 = -2307.5
$ws.Range("H74").Value2 = 1601.5
$ws.Range("I74").Value2 = 2120
$ws.Range("J74").Value2 = 1083
$ws.Range("K74").Value2 = 2120
$ws.Range("L74").Value2 = 1083
$ws.Range("M74").Value2 = -1246
$ws.Range("N74").Value2 = -2831
$ws.Range("H77").Value2 = 1601.5
$ws.Range("I77").Value2 = 2120
$ws.Range("J77").Value2 = 1083
$ws.Range("K77").Value2 = 10600
$ws.Range("L77").Value2 = 5415
$ws.Range("M77").Value2 = -6232
$ws.Range("N77").Value2 = -14151
$ws.Range("H96").Value2 = 45344
$ws.Range("J96").Value2 = 45344
$ws.Range("L96").Value2 = 45344
$ws.Range("N96").Value2 = -50836
$ws.Range("H125").Value2 = 49998.5
$ws.Range("J125").Value2 = 49998.5
$ws.Range("L125").Value2 = 49998.5
$ws.Range("N125").Value2 = -59838.5
$ws.Range("H132").Value2 = 2202.1924
$ws.Range("I132").Value2 = 1500.0256
$ws.Range("J132").Value2 = 4308.6924
$ws.Range("K132").Value2 = 4500.0768
$ws.Range("L132").Value2 = 12926.0772
$ws.Range("M132").Value2 = -1970.0768
$ws.Range("N132").Value2 = -17986.0772
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value2 = 160.5
$ws.Range("I4").Value2 = 480
$ws.Range("J4").Value2 = 54
$ws.Range("K4").Value2 = 480
$ws.Range("L4").Value2 = 54
$ws.Range("M4").Value2 = -365
$ws.Range("N4").Value2 = -284
$ws.Range("H95").Value2 = 47624
$ws.Range("J95").Value2 = 47624
$ws.Range("L95").Value2 = 47624
$ws.Range("N95").Value2 = -53116
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 1427.2
$ws.Range("I16").Value2 = 1370.3334
$ws.Range("J16").Value2 = 1512.5
$ws.Range("K16").Value2 = 1370.3334
$ws.Range("L16").Value2 = 1512.5
$ws.Range("M16").Value2 = -1083.3334
$ws.Range("N16").Value2 = -2086.5
$ws.Range("H113").Value2 = 1427.2
$ws.Range("I113").Value2 = 1370.3334
$ws.Range("J113").Value2 = 1512.5
$ws.Range("K113").Value2 = 1370.3334
$ws.Range("L113").Value2 = 1512.5
$ws.Range("M113").Value2 = 799.6666
$ws.Range("N113").Value2 = -5852.5
$ws.Range("H134").Value2 = 2038.6552
$ws.Range("I134").Value2 = 1652.0476
$ws.Range("J134").Value2 = 3053.5
$ws.Range("K134").Value2 = 4956.142800000001
$ws.Range("L134").Value2 = 9160.5
$ws.Range("M134").Value2 = -2421.142800000001
$ws.Range("N134").Value2 = -14230.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value2 = 4465071
$ws.Range("I2").Value2 = 1273.875
$ws.Range("J2").Value2 = 10416800
$ws.Range("K2").Value2 = 7643.25
$ws.Range("L2").Value2 = 62500800
$ws.Range("M2").Value2 = -7530.25
$ws.Range("N2").Value2 = -62501026
$ws.Range("H5").Value2 = 2244.5789
$ws.Range("I5").Value2 = 2538.5557
$ws.Range("K5").Value2 = 7615.6671
$ws.Range("M5").Value2 = -7503.6671
$ws.Range("H12").Value2 = 57.058823
$ws.Range("I12").Value2 = 16.5
$ws.Range("J12").Value2 = 69.53846
$ws.Range("K12").Value2 = 49.5
$ws.Range("L12").Value2 = 208.61538
$ws.Range("M12").Value2 = 123.5
$ws.Range("N12").Value2 = -554.61538
$ws.Range("H56").Value2 = 5423.4614
$ws.Range("I56").Value2 = 5423.4614
$ws.Range("K56").Value2 = 5423.4614
$ws.Range("M56").Value2 = -4893.4614
$ws.Range("H92").Value2 = 904.7143
$ws.Range("I92").Value2 = 600
$ws.Range("J92").Value2 = 955.5
$ws.Range("K92").Value2 = 1800
$ws.Range("L92").Value2 = 2866.5
$ws.Range("M92").Value2 = -552
$ws.Range("N92").Value2 = -5362.5
$ws.Range("H120").Value2 = 12287.454
$ws.Range("I120").Value2 = 9171.666999999999
$ws.Range("K120").Value2 = 27515.001
$ws.Range("M120").Value2 = -22677.001
$ws.Range("H122").Value2 = 760.3
$ws.Range("J122").Value2 = 1999.5
$ws.Range("L122").Value2 = 17995.5
$ws.Range("N122").Value2 = -22895.5
$ws.Range("H135").Value2 = 2244.5789
$ws.Range("I135").Value2 = 2538.5557
$ws.Range("K135").Value2 = 22847.0013
$ws.Range("M135").Value2 = -20312.0013
$ws.Range("H139").Value2 = 2259.1155
$ws.Range("I139").Value2 = 1874.1177
$ws.Range("J139").Value2 = 2986.3333
$ws.Range("K139").Value2 = 5622.3531
$ws.Range("L139").Value2 = 8958.999899999999
$ws.Range("M139").Value2 = -482.3531000000003
$ws.Range("N139").Value2 = -19238.9999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H117").Value2 = 0
$ws.Range("J117").Value2 = 0
$ws.Range("L117").Value2 = 0
$ws.Range("N117").ClearContents()
$ws.Range("H123").Value2 = 8890.177
$ws.Range("J123").Value2 = 8890.177
$ws.Range("L123").Value2 = 8890.177
$ws.Range("N123").Value2 = -13790.177
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value2 = 20658.666
$ws.Range("J76").Value2 = 20658.666
$ws.Range("L76").Value2 = 20658.666
$ws.Range("N76").Value2 = -21334.666
$ws.Range("H79").Value2 = 20658.666
$ws.Range("J79").Value2 = 20658.666
$ws.Range("L79").Value2 = 20658.666
$ws.Range("N79").Value2 = -22998.666
$ws.Range("H100").Value2 = 5930.385
$ws.Range("I100").Value2 = 8656.571
$ws.Range("J100").Value2 = 2749.8333
$ws.Range("K100").Value2 = 8656.571
$ws.Range("L100").Value2 = 2749.8333
$ws.Range("M100").Value2 = -8115.571
$ws.Range("N100").Value2 = -3831.8333
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value2 = 34990
$ws.Range("J82").Value2 = 34990
$ws.Range("L82").Value2 = 34990
$ws.Range("N82").Value2 = -35756
$ws.Range("H85").Value2 = 34990
$ws.Range("J85").Value2 = 34990
$ws.Range("L85").Value2 = 34990
$ws.Range("N85").Value2 = -37642
$ws.Range("H101").Value2 = 166671790
$ws.Range("J101").Value2 = 166671790
$ws.Range("L101").Value2 = 166671790
$ws.Range("N101").Value2 = -166678280
$ws.Range("H132").Value2 = 2611.0312
$ws.Range("I132").Value2 = 1679.4117
$ws.Range("J132").Value2 = 3666.8667
$ws.Range("K132").Value2 = 5038.2351
$ws.Range("L132").Value2 = 11000.6001
$ws.Range("M132").Value2 = -2508.2351
$ws.Range("N132").Value2 = -16060.6001
$ws.Range("H136").Value2 = 1652.5
$ws.Range("I136").Value2 = 1326.8518
$ws.Range("K136").Value2 = 3980.5554
$ws.Range("M136").Value2 = -1430.5554
